$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Policy Number -> N/A ---
$ws.Range("B5").Value = "N/A"

# --- Insert 3 new data rows ---
# 1) New row at (old) position 13 -> pushes Standard pregnancy test (old 13) to 14,
#    Prostatectomy (old 14) to 15, seasonique (old 15) to 16, Total (old 16) to 17.
$ws.Rows.Item(13).Insert()
# 2) Two more rows after the (now shifted) "Standard pregnancy test" row (14),
#    i.e. at rows 15 and 16 -> pushes Prostatectomy to 17, seasonique to 18, Total to 19.
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(16).Insert()

# Rows inherit the format of the row above on insert. Row 13 correctly inherits the
# "Fraud" (red) look from row 12 (Colonoscopy). Rows 15-16 inherit "Legitimate" (green)
# from row 14 (Standard pregnancy test) - fix them to "Risk" (orange) by copying the
# format from an existing Risk row (row 9) so the same style index is reused.
$ws.Range("A9:E9").Copy()
$ws.Range("A15:E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 13: Spirometry (procedure) - Fraud
$ws.Range("A13").Value = "Spirometry (procedure)"
$ws.Range("B13").Value = 15000
$ws.Range("C13").Value = 7786.47
$ws.Range("D13").Value = 7213.53
$ws.Range("E13").Value = "Fraud"

# Row 15: Throat culture (procedure) - Risk
$ws.Range("A15").Value = "Throat culture (procedure)"
$ws.Range("B15").Value = 2300
$ws.Range("C15").Value = 2020.43
$ws.Range("D15").Value = 279.57
$ws.Range("E15").Value = "Risk"

# Row 16: Upper arm X-ray - Risk
$ws.Range("A16").Value = "Upper arm X-ray"
$ws.Range("B16").Value = 1500
$ws.Range("C16").Value = 431.4
$ws.Range("D16").Value = 1068.6
$ws.Range("E16").Value = "Risk"

# --- Update Total Invoice Amount (now on row 19) ---
$ws.Range("B19").Value = 66433.10000000001

# --- Append "Overall Status" row 20, styled like a Fraud row (copy from row 12) ---
$ws.Range("A12:E12").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A20").Value = "Overall Status"
$ws.Range("E20").Value = "Fraud"
